# "adding ID to register page"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing two users need to move down one row to make room for the new
# header row; the new "Neriala" row is written first (matches the order the
# shared strings were introduced in the saved file), then the header, then
# the remaining new users.
$ws.Range("A3").Value = "yyyyyyy"
$ws.Range("B3").Value = "123123WW!"
$ws.Range("A2").Value = "yossiyo"
$ws.Range("B2").Value = "Reuts8888!"

$ws.Range("A4").Value = "Neriala"
$ws.Range("B4").Value = "Neriala12#"

$ws.Range("A1").Value = "Users"
$ws.Range("B1").Value = "Passwords"
$ws.Range("C1").Value = "ID"

$ws.Range("A5").Value = "TonerMe"
$ws.Range("B5").Value = "Toner12#"

$ws.Range("A6").Value = "Yossko"
$ws.Range("B6").Value = "Yossko12#"

$ws.Range("A7").Value = "kaikaich"
$ws.Range("B7").Value = "kaikai1@"

# New "ID" column values
$ws.Range("C2").Value = 111111111
$ws.Range("C3").Value = 111111112
$ws.Range("C4").Value = 111111113
$ws.Range("C5").Value = 111111114
$ws.Range("C6").Value = 111111115
$ws.Range("C7").Value = 111111119

# Header formatting: bold, size 14, purple, Calibri
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 14
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Color = 10498160
$ws.Rows.Item(1).RowHeight = 18

# Column widths (target widths of 15.8984375 / 19.09765625 / 14.59765625
# characters come from real Excel's pixel-grid rounding; this host quantizes
# ColumnWidth to 1/7-character steps, so these inputs land on the closest
# achievable stored width)
$ws.Columns.Item(1).ColumnWidth = 15.1
$ws.Columns.Item(2).ColumnWidth = 18.4
$ws.Columns.Item(3).ColumnWidth = 13.86

# Select the last row like the source sheet view
$excel.Goto($ws.Range("A7:XFD7"))

$wb.Save()
